# Apply the commit: "changed gfx size, enable all keys, correct LP phases (resorted)"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- View changes: zoom to 100%, selection moved to C10 ---
$win = $ws.Application.ActiveWindow
$win.Zoom = 100
$ws.Range("C10").Select()

# --- Correct LP phases (resorted): update columns A, B, D, F for the affected rows ---
$ws.Range("A5").Value = 39
$ws.Range("B5").Value = "si"
$ws.Range("D5").Value = "ko"
$ws.Range("F5").Value = "me"

$ws.Range("A12").Value = 6
$ws.Range("B12").Value = "go"
$ws.Range("D12").Value = "fo"
$ws.Range("F12").Value = "ta"

$ws.Range("A13").Value = 45
$ws.Range("D13").Value = "no"

$ws.Range("A14").Value = 5
$ws.Range("B14").Value = "go"
$ws.Range("D14").Value = "fa"
$ws.Range("F14").Value = "ta"

$ws.Range("A15").Value = 49
$ws.Range("D15").Value = "pa"

$ws.Range("A16").Value = 52
$ws.Range("B16").Value = "si"
$ws.Range("D16").Value = "to"
$ws.Range("F16").Value = "me"

$ws.Range("A17").Value = 3
$ws.Range("D17").Value = "do"

$ws.Range("A18").Value = 21
$ws.Range("D18").Value = "pa"

$ws.Range("A26").Value = 34
$ws.Range("D26").Value = "fo"

$ws.Range("A27").Value = 1
$ws.Range("B27").Value = "go"
$ws.Range("D27").Value = "ba"
$ws.Range("F27").Value = "ta"

$ws.Range("A28").Value = 48
$ws.Range("B28").Value = "si"
$ws.Range("D28").Value = "su"
$ws.Range("F28").Value = "me"

$ws.Range("A29").Value = 25
$ws.Range("B29").Value = "go"
$ws.Range("D29").Value = "tu"
$ws.Range("F29").Value = "ta"

$ws.Range("A30").Value = 35
$ws.Range("B30").Value = "si"
$ws.Range("D30").Value = "fu"
$ws.Range("F30").Value = "me"

$ws.Range("A31").Value = 46
$ws.Range("D31").Value = "nu"

$ws.Range("A32").Value = 27
$ws.Range("B32").Value = "go"
$ws.Range("D32").Value = "za"
$ws.Range("F32").Value = "ta"

$ws.Range("A33").Value = 19
$ws.Range("D33").Value = "so"

$ws.Range("A34").Value = 37
$ws.Range("B34").Value = "si"
$ws.Range("D34").Value = "hi"
$ws.Range("F34").Value = "me"

$ws.Range("A35").Value = 10
$ws.Range("D35").Value = "ka"

$ws.Range("A36").Value = 50
$ws.Range("D36").Value = "po"

$ws.Range("A37").Value = 26
$ws.Range("D37").Value = "ti"

$ws.Range("A38").Value = 32
$ws.Range("D38").Value = "du"

$ws.Range("A39").Value = 13
$ws.Range("D39").Value = "la"

$ws.Range("A40").Value = 33
$ws.Range("D40").Value = "fa"
